$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-10-10 Thursday" "2024-10-11 Friday"
Replace-Text "75÷4=" "89÷4="
Replace-Text "45÷5=" "61÷4="
Replace-Text "10÷5=" "30÷9="
Replace-Text "55÷6=" "34÷6="
Replace-Text "14÷8=" "79÷6="
Replace-Text "40÷5=" "89÷8="
Replace-Text "52÷7=" "76÷8="
Replace-Text "30÷5=" "21÷5="
Replace-Text "16÷4=" "52÷8="
Replace-Text "59÷9=" "60÷8="
Replace-Text "74÷8=" "17÷9="
Replace-Text "64÷4=" "66÷2="
Replace-Text "38÷9=" "75÷2="
Replace-Text "39÷3=" "88÷9="
Replace-Text "96÷4=" "68÷6="
Replace-Text "66÷4=" "70÷7="
Replace-Text "90÷7=" "85÷9="
Replace-Text "10÷7=" "46÷9="
Replace-Text "65÷3=" "77÷5="
Replace-Text "18÷2=" "52÷6="
Replace-Text "82÷9=" "28÷8="
Replace-Text "20÷4=" "66÷7="
Replace-Text "59÷6=" "55÷2="
Replace-Text "48÷9=" "94÷9="
Replace-Text "25÷5=" "65÷3="

Write-Output "Done applying replacements"
